$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# CONTENT
# ===========================================================================

# ---------------------------------------------------------------------------
# 1) Row 4 ("Number of disability persons") becomes the new
#    "family with disabilities Persons " row with brand new figures.
#    (Written before the title so the shared-string table grows in the same
#    order the source workbook shows: family-row, disabilities-row, title.)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 1227
$ws.Range("C4").Value = 1173
$ws.Range("D4").Value = 1116
$ws.Range("E4").Value = 1148
$ws.Range("F4").Value = 1144
$ws.Range("G4").Value = 1183
$ws.Range("H4").Value = 1178
$ws.Range("I4").Value = 1138
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row 5 carrying the second data series
#    "disabilities Persons " (this pushes the old source row from 5 -> 6).
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 1395
$ws.Range("C5").Value = 1355
$ws.Range("D5").Value = 1291
$ws.Range("E5").Value = 1327
$ws.Range("F5").Value = 1328
$ws.Range("G5").Value = 1386
$ws.Range("H5").Value = 1393
$ws.Range("I5").Value = 1356
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 3) Row 1 title text is rewritten and merged across A1:I1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Lagodekhi Municipality"
$ws.Range("A1:I1").Merge()
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 4) The old source row (now row 6) keeps its "Source: ..." text - only the
#    merge range shifts down from A5:H5 to A6:H6 (handled automatically by
#    the row insert above). It does get a slightly taller row height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 27.75

# ===========================================================================
# FORMATTING
# ===========================================================================

# --- Column A width -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.95

# --- Row 1: new title, bold Arial 11, centered, wrapped --------------------
$r1 = $ws.Range("A1:I1")
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true
$r1.Interior.Pattern = -4142
$r1.Borders.LineStyle = -4142

# --- Row 2: "(End of year, persons)" Arial 10 --------------------------
$r2 = $ws.Range("A2")
$r2.Font.Name = "Arial"
$r2.Font.Size = 10
$r2.Font.ColorIndex = 1
$r2.Interior.Pattern = -4142
$r2.Borders.LineStyle = -4142

# --- Row 3: A3 blank cell, Sylfaen 11 with top border -----------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Sylfaen"
$a3.Font.Size = 11
$a3.Font.ColorIndex = 1
$a3.Interior.Pattern = -4142
$a3.Borders.Item(8).LineStyle = 1
$a3.Borders.Item(8).Weight = 2

# Row 3: years header B3:I3 - Arial 10 black, centered, grey fill, top+bottom border
$yr = $ws.Range("B3:I3")
$yr.Font.Name = "Arial"
$yr.Font.Size = 10
$yr.Font.ColorIndex = 1
$yr.HorizontalAlignment = -4108
$yr.VerticalAlignment = -4108
$yr.WrapText = $true
$yr.Interior.Pattern = 1
$yr.Interior.PatternColorIndex = -4142
$yr.Interior.ThemeColor = 2
$yr.Interior.TintAndShade = 0
$yr.Borders.Item(8).LineStyle = 1
$yr.Borders.Item(8).Weight = 2
$yr.Borders.Item(9).LineStyle = 1
$yr.Borders.Item(9).Weight = 2

# --- Row 4: "family with disabilities Persons " label + values -------------
$a4 = $ws.Range("A4")
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.ColorIndex = 1
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true
$a4.Interior.Pattern = 1
$a4.Interior.PatternColorIndex = -4142
$a4.Interior.ThemeColor = 2
$a4.Interior.TintAndShade = 0
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = 2

$v4 = $ws.Range("B4:I4")
$v4.Font.Name = "Arial"
$v4.Font.Size = 10
$v4.Font.ColorIndex = 1
$v4.NumberFormat = "#\ ##0"
$v4.Interior.Pattern = 1
$v4.Interior.PatternColorIndex = -4142
$v4.Interior.ThemeColor = 2
$v4.Interior.TintAndShade = 0
$v4.Borders.LineStyle = -4142

# --- Row 5: "disabilities Persons " label + values --------------------------
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.ColorIndex = 1
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$a5.Interior.Pattern = 1
$a5.Interior.PatternColorIndex = -4142
$a5.Interior.ThemeColor = 2
$a5.Interior.TintAndShade = 0
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = 2

$v5 = $ws.Range("B5:H5")
$v5.Font.Name = "Arial"
$v5.Font.Size = 10
$v5.Font.ColorIndex = 1
$v5.NumberFormat = "#\ ##0"
$v5.Interior.Pattern = 1
$v5.Interior.PatternColorIndex = -4142
$v5.Interior.ThemeColor = 2
$v5.Interior.TintAndShade = 0
$v5.Borders.LineStyle = -4142

$i5 = $ws.Range("I5")
$i5.Font.Name = "Arial"
$i5.Font.Size = 10
$i5.Font.ColorIndex = 1
$i5.NumberFormat = "#\ ##0"
$i5.Interior.Pattern = 1
$i5.Interior.PatternColorIndex = -4142
$i5.Interior.ThemeColor = 2
$i5.Interior.TintAndShade = 0
$i5.Borders.Item(9).LineStyle = 1
$i5.Borders.Item(9).Weight = 2

# --- Row 6: source text row, Arial 9, left/center/wrap ----------------------
$a6 = $ws.Range("A6")
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Font.ColorIndex = 1
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true
$a6.Interior.Pattern = 1
$a6.Interior.PatternColorIndex = -4142
$a6.Interior.ThemeColor = 2
$a6.Interior.TintAndShade = 0
$a6.Borders.LineStyle = -4142

$v6 = $ws.Range("B6:H6")
$v6.Font.Name = "Arial"
$v6.Font.Size = 9
$v6.Font.ColorIndex = 1
$v6.HorizontalAlignment = -4131
$v6.VerticalAlignment = -4108
$v6.WrapText = $true
$v6.Interior.Pattern = 1
$v6.Interior.PatternColorIndex = -4142
$v6.Interior.ThemeColor = 2
$v6.Interior.TintAndShade = 0
$v6.Borders.Item(8).LineStyle = 1
$v6.Borders.Item(8).Weight = 2

# --- Selection / active cell moves to the new title block -------------------
$ws.Range("A1:I1").Select()

$wb.Save()
